$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 635 (the "「何事も永遠には続かない」" entry) entirely.
# This shifts all subsequent rows up by one, matching the target diff.
$ws.Rows.Item(635).Delete()
